$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "chr" in F1, matching style of existing header cells
$ws.Range("F1").Value = "chr"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill in chromosome values for rows 2-13
$values = @(10, 9, 6, 6, 8, 6, 6, 6, 6, 17, 6, 14)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
